# Scheduled runner update: refresh market-derived price/profit columns (H:N)
# across the Leve profit-tracking sheets with the latest Universalis snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 44124.883
$ws.Range("J70").Value = 64507.176
$ws.Range("L70").Value = 193521.528
$ws.Range("N70").Value = -194061.528
# Row 73
$ws.Range("H73").Value = 44124.883
$ws.Range("J73").Value = 64507.176
$ws.Range("L73").Value = 193521.528
$ws.Range("N73").Value = -195393.528
# Row 106
$ws.Range("H106").Value = 2344.25
$ws.Range("I106").Value = 2344.25
$ws.Range("K106").Value = 2344.25
$ws.Range("M106").Value = -1713.25
# Row 112
$ws.Range("H112").Value = 1411.1666
$ws.Range("J112").Value = 1443.4
$ws.Range("L112").Value = 4330.200000000001
$ws.Range("N112").Value = -6546.200000000001
# Row 135
$ws.Range("H135").Value = 2420.1
$ws.Range("J135").Value = 4985.5
$ws.Range("L135").Value = 44869.5
$ws.Range("N135").Value = -49939.5
# Row 138
$ws.Range("H138").Value = 4030.3774
$ws.Range("I138").Value = 3485
$ws.Range("K138").Value = 10455
$ws.Range("M138").Value = -5315

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6587.5835
$ws.Range("I32").Value = 6587.5835
$ws.Range("K32").Value = 6587.5835
$ws.Range("M32").Value = -6300.5835
# Row 62
$ws.Range("H62").Value = 31993
$ws.Range("J62").Value = 31993
$ws.Range("L62").Value = 31993
$ws.Range("N62").Value = -33241
# Row 63
$ws.Range("H63").Value = 2689.8572
$ws.Range("I63").Value = 2582.6667
$ws.Range("J63").Value = 3333
$ws.Range("K63").Value = 2582.6667
$ws.Range("L63").Value = 3333
$ws.Range("M63").Value = -1896.6667
$ws.Range("N63").Value = -4705
# Row 65
$ws.Range("H65").Value = 31993
$ws.Range("J65").Value = 31993
$ws.Range("L65").Value = 95979
$ws.Range("N65").Value = -102219
# Row 66
$ws.Range("H66").Value = 2689.8572
$ws.Range("I66").Value = 2582.6667
$ws.Range("J66").Value = 3333
$ws.Range("K66").Value = 12913.3335
$ws.Range("L66").Value = 16665
$ws.Range("M66").Value = -9481.333500000001
$ws.Range("N66").Value = -23529
# Row 74
$ws.Range("H74").Value = 2185.0952
$ws.Range("I74").Value = 2050.5278
$ws.Range("K74").Value = 2050.5278
$ws.Range("M74").Value = -1176.5278
# Row 77
$ws.Range("H77").Value = 2185.0952
$ws.Range("I77").Value = 2050.5278
$ws.Range("K77").Value = 10252.639
$ws.Range("M77").Value = -5884.638999999999
# Row 97
$ws.Range("H97").Value = 1250
$ws.Range("I97").Value = 951.5714
$ws.Range("J97").Value = 2294.5
$ws.Range("K97").Value = 951.5714
$ws.Range("L97").Value = 2294.5
$ws.Range("M97").Value = -455.5714
$ws.Range("N97").Value = -3286.5

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 864.6
$ws.Range("J94").Value = 798
$ws.Range("L94").Value = 798
$ws.Range("N94").Value = -1700

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2930.111
$ws.Range("I16").Value = 3232
$ws.Range("J16").Value = 2326.3333
$ws.Range("K16").Value = 3232
$ws.Range("L16").Value = 2326.3333
$ws.Range("M16").Value = -2945
$ws.Range("N16").Value = -2900.3333
# Row 50
$ws.Range("H50").Value = 9998.888999999999
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 9998.888999999999
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 9998.888999999999
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -11248.889
# Row 62
$ws.Range("H62").Value = 3464.5
$ws.Range("I62").Value = 2929.6667
$ws.Range("K62").Value = 2929.6667
$ws.Range("M62").Value = -2305.6667
# Row 65
$ws.Range("H65").Value = 3464.5
$ws.Range("I65").Value = 2929.6667
$ws.Range("K65").Value = 14648.3335
$ws.Range("M65").Value = -11528.3335
# Row 113
$ws.Range("H113").Value = 2930.111
$ws.Range("I113").Value = 3232
$ws.Range("J113").Value = 2326.3333
$ws.Range("K113").Value = 3232
$ws.Range("L113").Value = 2326.3333
$ws.Range("M113").Value = -1062
$ws.Range("N113").Value = -6666.3333
# Row 125
$ws.Range("H125").Value = 54499.5
$ws.Range("J125").Value = 54499.5
$ws.Range("L125").Value = 54499.5
$ws.Range("N125").Value = -59419.5

$ws = $wb.Worksheets.Item("CUL")
# Row 75
$ws.Range("H75").Value = 4992.2
$ws.Range("I75").Value = 4992
$ws.Range("J75").Value = 4992.25
$ws.Range("K75").Value = 14976
$ws.Range("L75").Value = 14976.75
$ws.Range("M75").Value = -13978
$ws.Range("N75").Value = -16972.75
# Row 78
$ws.Range("H78").Value = 4992.2
$ws.Range("I78").Value = 4992
$ws.Range("J78").Value = 4992.25
$ws.Range("K78").Value = 44928
$ws.Range("L78").Value = 44930.25
$ws.Range("M78").Value = -39936
$ws.Range("N78").Value = -54914.25

$ws = $wb.Worksheets.Item("GSM")
# Row 22
$ws.Range("H22").Value = 337.42856
$ws.Range("I22").Value = 347.83334
$ws.Range("K22").Value = 347.83334
$ws.Range("M22").Value = 181.16666
# Row 96
$ws.Range("H96").Value = 19970
$ws.Range("J96").Value = 19970
$ws.Range("L96").Value = 19970
$ws.Range("N96").Value = -25462
# Row 97
$ws.Range("H97").Value = 366.4
$ws.Range("J97").Value = 297
$ws.Range("L97").Value = 297
$ws.Range("N97").Value = -1289
# Row 122
$ws.Range("H122").Value = 2672
$ws.Range("I122").Value = 1477
$ws.Range("K122").Value = 4431
$ws.Range("M122").Value = -1981

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2298
$ws.Range("J22").Value = 2667.2222
$ws.Range("L22").Value = 2667.2222
$ws.Range("N22").Value = -3257.2222
# Row 27
$ws.Range("H27").Value = 2298
$ws.Range("J27").Value = 2667.2222
$ws.Range("L27").Value = 2667.2222
$ws.Range("N27").Value = -2881.2222
# Row 46
$ws.Range("H46").Value = 27823.426
$ws.Range("I46").Value = 57774.168
$ws.Range("J46").Value = 3318.2727
$ws.Range("K46").Value = 57774.168
$ws.Range("L46").Value = 3318.2727
$ws.Range("M46").Value = -57586.168
$ws.Range("N46").Value = -3694.2727
# Row 55
$ws.Range("H55").Value = 1497.8125
$ws.Range("I55").Value = 1305.4445
$ws.Range("J55").Value = 1745.1428
$ws.Range("K55").Value = 1305.4445
$ws.Range("L55").Value = 1745.1428
$ws.Range("M55").Value = -1132.4445
$ws.Range("N55").Value = -2091.1428
# Row 93
$ws.Range("H93").Value = 491.63635
$ws.Range("I93").Value = 434.2857
$ws.Range("J93").Value = 592
$ws.Range("K93").Value = 434.2857
$ws.Range("L93").Value = 592
$ws.Range("M93").Value = 813.7143
$ws.Range("N93").Value = -3088
